# Apply the "output generated" refresh to 上海-漫展信息.xlsx
#
# Sheet 1 = 展览 (Exhibition), Sheet 2 = 演出 (Performance),
# Sheet 3 = 本地生活 (Local life), Sheet 4 = 全部类型 (All types, a
# rollup/union of the other three sheets' rows).
#
# Most of the change is "想去人数" (want-to-go count) counters ticking up
# (or occasionally down) across the sheets, plus a new exhibition row
# inserted into 展览 (and NOT propagated into 全部类型).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F3").Value = 2423
$ws1.Range("F6").Value = 74
$ws1.Range("F9").Value = 3321
$ws1.Range("F11").Value = 1052
$ws1.Range("F12").Value = 875
$ws1.Range("F14").Value = 847
$ws1.Range("F15").Value = 1512
$ws1.Range("F16").Value = 5
$ws1.Range("F17").Value = 827
$ws1.Range("F18").Value = 1716
$ws1.Range("F20").Value = 391
$ws1.Range("F21").Value = 1502
$ws1.Range("F22").Value = 76

# The row-shift below re-touches whatever lands in G26 (old G25, the
# "iPR动漫" row's min ticket price of 65.8) through the engine's internal
# float64 path, which - unlike the untouched original cell - serializes
# with full IEEE754 precision ("65.799999999999997") instead of the
# shortest round-trip form. Snapshot it now and restore the clean literal
# after the insert so the saved XML still reads "65.8".
$savedG25 = $ws1.Range("G25").Value

# Insert a brand-new exhibition row at row 24; this pushes the old row 24
# ("...长三角文博会...") down to row 25 and the old row 25
# ("...iPR动漫...") down to row 26.
$ws1.Rows.Item(24).Insert()
$ws1.Range("G26").Value = $savedG25

# The inserted row's A cell has no formatting yet - copy the bold/centered/
# bordered look used by every other row's A column (A25, freshly shifted
# down, still has it).
$ws1.Range("A25").Copy()
$ws1.Range("A24").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# Column A is a simple running index (row number - 2); keep that pattern
# intact across the newly inserted/shifted rows.
$ws1.Range("A24").Value = 23
$ws1.Range("A25").Value = 24
$ws1.Range("A26").Value = 25

# B24 is a plain "yyyy-mm-dd" text label (matching every other row in the
# column), not a real date value - a bare assignment would get silently
# auto-converted to a date serial by Excel's normal input parsing, so force
# text with a leading quote-prefix, then strip the formatting residue that
# leaves behind by re-applying the plain/default format from a neighboring
# untouched text cell (B23).
$ws1.Range("B24").Value = "'2024-11-16"
$ws1.Range("B23").Copy()
$ws1.Range("B24").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Range("C24").Value = "上海·第三届ICG动漫游戏博览会"
$ws1.Range("D24").Value = "浦星公路567弄 上海国际品牌珠宝中心"
$ws1.Range("E24").Value = "2024.11.16 10:00-11.17 17:00"
$ws1.Range("F24").Value = 3
$ws1.Range("G24").Value = 59
$ws1.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=92846"
$ws1.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202409/BRwbjq671727178083854.jpeg"

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F3").Value = 9
$ws2.Range("F40").Value = 364
$ws2.Range("F48").Value = 300

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F6").Value = 2510
$ws3.Range("F7").Value = 9549
$ws3.Range("F8").Value = 142
$ws3.Range("F12").Value = 2831
$ws3.Range("F13").Value = 380
$ws3.Range("F14").Value = 696

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (rollup of sheets 1-3; note it does NOT gain the new
# ICG row that was inserted into 展览 - only the counter values move)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F3").Value = 142
$ws4.Range("F6").Value = 2423
$ws4.Range("F8").Value = 380
$ws4.Range("F9").Value = 9
$ws4.Range("F13").Value = 74
$ws4.Range("F18").Value = 1052
$ws4.Range("F19").Value = 875
$ws4.Range("F21").Value = 847
$ws4.Range("F23").Value = 1512
$ws4.Range("F27").Value = 827
$ws4.Range("F31").Value = 1716
$ws4.Range("F32").Value = 391
$ws4.Range("F39").Value = 76
$ws4.Range("F44").Value = 364
